# Update crypto price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.762.36'
$ws.Range('E2').Value = '  +3.51%  '
$ws.Range('D3').Value = '2.257.16'
$ws.Range('E3').Value = '  +3.29%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '''253.73'
$ws.Range('E5').Value = '  -0.62%  '
$ws.Range('D6').Value = '''0.626'
$ws.Range('E6').Value = '  -0.19%  '
$ws.Range('D7').Value = '''71.95'
$ws.Range('E7').Value = '  +5.64%  '
$ws.Range('D9').Value = '''0.644'
$ws.Range('E9').Value = '  +11.92%  '
$ws.Range('D10').Value = '''41.21'
$ws.Range('E10').Value = '  +9.43%  '
$ws.Range('D11').Value = '''59.68'
$ws.Range('E11').Value = '  +1.03%  '
$ws.Range('D12').Value = '''0.0967'
$ws.Range('E12').Value = '  +3.46%  '
$ws.Range('D13').Value = '''7.38'
$ws.Range('E13').Value = '  +3.71%  '
$ws.Range('D14').Value = '''0.104'
$ws.Range('E14').Value = '  -0.28%  '
$ws.Range('D15').Value = '2.595.41'
$ws.Range('E15').Value = '  +3.71%  '
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').Value = '''0.884'
$ws.Range('E16').Value = '  +1.49%  '
$ws.Range('B17').Value = 'Chainlink'
$ws.Range('C17').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D17').Value = '''14.76'
$ws.Range('E17').Value = '  +2.16%  '
$ws.Range('D18').Value = '2.254.36'
$ws.Range('E18').Value = '  +1.86%  '
$ws.Range('D19').Value = '42.761.97'
$ws.Range('E19').Value = '  +3.78%  '
$ws.Range('D20').Value = '0.0₃0982'
$ws.Range('E20').Value = '  +2.96%  '
$ws.Range('D21').Value = '''6.26'
$ws.Range('E21').Value = '  +1.41%  '
$ws.Range('D22').Value = '''73.23'
$ws.Range('E22').Value = '  +1.92%  '
$ws.Range('D23').Value = '''234.15'
$ws.Range('E23').Value = '  +0.87%  '
$ws.Range('E24').Value = '  +3.91%  '
$ws.Range('E25').Value = '  +1.09%  '
$ws.Range('D26').Value = '''11.73'
$ws.Range('E26').Value = '  +0.35%  '
$ws.Range('D27').Value = '''0.999'
$ws.Range('E27').Value = '  -0.01%  '
$ws.Range('D28').Value = '''2.47'
$ws.Range('E28').Value = '  -2.74%  '
$ws.Range('D29').Value = '''3.66'
$ws.Range('E29').Value = '  -2.28%  '
$ws.Range('D30').Value = '''2.15'
$ws.Range('E30').Value = '  -0.89%  '
$ws.Range('D31').Value = '''167.67'
$ws.Range('E31').Value = '  -0.70%  '
$ws.Range('D32').Value = '''21.02'
$ws.Range('E32').Value = '  +1.77%  '
$ws.Range('D33').Value = '''0.132'
$ws.Range('E33').Value = '  +12.48%  '
$ws.Range('D34').Value = '''6.13'
$ws.Range('E34').Value = '  +12.13%  '
$ws.Range('D35').Value = '''0.0788'
$ws.Range('E35').Value = '  +4.85%  '
$ws.Range('D36').Value = '''0.125'
$ws.Range('E36').Value = '  +1.39%  '
$ws.Range('D37').Value = '''28.88'
$ws.Range('E37').Value = '  +9.23%  '
$ws.Range('D38').Value = '''4.71'
$ws.Range('E38').Value = '  +1.86%  '
$ws.Range('D39').Value = '''4.17'
$ws.Range('E39').Value = '  +0.61%  '
$ws.Range('D40').Value = '''0.0321'
$ws.Range('E40').Value = '  +6.70%  '
$ws.Range('D41').Value = '''2.30'
$ws.Range('E41').Value = '  +4.56%  '
$ws.Range('D42').Value = '''6.12'
$ws.Range('E42').Value = '  +8.07%  '
$ws.Range('D43').Value = '''12.47'
$ws.Range('E43').Value = '  -0.71%  '
$ws.Range('E44').Value = '  +0.71%  '
$ws.Range('E45').Value = '  -3.82%  '
$ws.Range('E46').Value = '  +0.55%  '
$ws.Range('D47').Value = '''8.95'
$ws.Range('E47').Value = '  +3.63%  '
$ws.Range('E48').Value = '  +1.15%  '
$ws.Range('E49').Value = '  +4.70%  '
$ws.Range('E50').Value = '  -0.40%  '
$ws.Range('D51').Value = '''4.43'
$ws.Range('E51').Value = '  +4.22%  '
